# Update scripts with new TPM data: refresh existing row 2 values and
# append a new row 3 for the Il12b -> Il23r ligand-receptor pair table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Sending/Target cluster labels and refresh all computed values ---
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Il12b"
$ws.Range("C2").Value = "Il23r"
$ws.Range("D2").Value = "MuSCs"

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.284858666666667
$ws.Range("H2").Value = 3.854576
$ws.Range("I2").Value = 0.6825120125588942
$ws.Range("J2").Value = 0.6825120125588942
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04298466666666667
$ws.Range("N2").Value = 0.128954
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.05522922150044446
$ws.Range("R2").Value = 0.4970629935040001
$ws.Range("S2").Value = 0.6825120125588942
$ws.Range("T2").Value = 0.6825120125588942

# --- Row 3: brand new row with its own Sending/Target cluster labels and values ---
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Il12b"
$ws.Range("C3").Value = "Il23r"
$ws.Range("D3").Value = "MuSCs"

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.597685
$ws.Range("H3").Value = 1.793055
$ws.Range("I3").Value = 0.3174879874411058
$ws.Range("J3").Value = 0.3174879874411058
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04298466666666667
$ws.Range("N3").Value = 0.128954
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.02569129049666667
$ws.Range("R3").Value = 0.23122161447
$ws.Range("S3").Value = 0.3174879874411058
$ws.Range("T3").Value = 0.3174879874411058
